# Generate Report for Handoff
# Replaces the "handed back" report rows with a freshly generated
# "ready for handoff" report: new file GUIDs, new status/timestamps,
# and drops the now-unused Latest Target File / Latest Handback File
# columns (F/G) on the per-locale sheets.

$wb = $excel.ActiveWorkbook

# ---- identifiers that changed ----------------------------------------
# old: 11f943d4-e580-437c-8da3-f377f2e05c15 / bb898a2c-ec1c-4f15-987d-4e981b05dd91
$newGuid1 = "56426c72-e1a0-49f7-bb3a-df81ef7ff23f"
$newGuid2 = "ffffcf2e933f-cbf8-4297-a5f3-fbd5915e14d0"

$newHash  = "8d2f12b505b9981efa8ba5b178d3047ae1a44b04"

$status      = "Ready for handoff"
$latestDate  = "2016-50-18 22:50:35"

$zhFile      = "$newGuid1.$newHash.zh-cn.xlf"
$deFile      = "$newGuid1.$newHash.de-de.xlf"

$zhHandoffDt = "2016-03-18 22:50:32"
$deHandoffDt = "2016-03-18 22:50:35"
$handbackDt  = "0001-01-01 00:00:00"

$md1 = "$newGuid1.md"
$md2 = "$newGuid2.md"

function Base-Url($guid) {
    return "https://github.com/OpenLocalizationTest/oltest/blob/da13284d45f18ef3a992b8727d432a18cc4e8b4e/e2e/$guid.md"
}

# =========================================================================
# Sheet "Overview"
# =========================================================================
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A1").Hyperlinks.Delete() | Out-Null

$ws1.Range("B2").Value = $status
$ws1.Range("C2").Value = $status
$ws1.Range("D2").Value = $latestDate

$ws1.Range("B3").Value = $status
$ws1.Range("C3").Value = $status
$ws1.Range("D3").Value = $latestDate

$ws1.Hyperlinks.Add($ws1.Range("A2"), (Base-Url $newGuid1), "", "", $md1) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), (Base-Url $newGuid2), "", "", $md2) | Out-Null

# =========================================================================
# Sheet "zh-cn"
# =========================================================================
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A1").Hyperlinks.Delete() | Out-Null

$ws2.Range("C2").Value = $status
$ws2.Range("D2").Value = $zhFile
$ws2.Range("E2").Value = $zhHandoffDt
$ws2.Range("H2").Value = $handbackDt

$ws2.Range("C3").Value = $status
$ws2.Range("D3").Value = $zhFile
$ws2.Range("E3").Value = $zhHandoffDt
$ws2.Range("H3").Value = $handbackDt

# Latest Target File / Latest Handback File columns no longer populated
$ws2.Range("F2:G3").Clear() | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A2"), (Base-Url $newGuid1), "", "", $md1) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B2"), (Base-Url $newGuid1), "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c7cae3eadf7f18b3c548cd26d9555be48b7c6f7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhFile", "", "", $zhFile) | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A3"), (Base-Url $newGuid2), "", "", $md2) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B3"), (Base-Url $newGuid2), "", "", ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c7cae3eadf7f18b3c548cd26d9555be48b7c6f7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhFile", "", "", $zhFile) | Out-Null

# =========================================================================
# Sheet "de-de"
# =========================================================================
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A1").Hyperlinks.Delete() | Out-Null

$ws3.Range("C2").Value = $status
$ws3.Range("D2").Value = $deFile
$ws3.Range("E2").Value = $deHandoffDt
$ws3.Range("H2").Value = $handbackDt

$ws3.Range("C3").Value = $status
$ws3.Range("D3").Value = $deFile
$ws3.Range("E3").Value = $deHandoffDt
$ws3.Range("H3").Value = $handbackDt

# Latest Target File / Latest Handback File columns no longer populated
$ws3.Range("F2:G3").Clear() | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A2"), (Base-Url $newGuid1), "", "", $md1) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B2"), (Base-Url $newGuid1), "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7adf40773f59f8bce68bcace7ba22d63e77bcecb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deFile", "", "", $deFile) | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A3"), (Base-Url $newGuid2), "", "", $md2) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B3"), (Base-Url $newGuid2), "", "", ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7adf40773f59f8bce68bcace7ba22d63e77bcecb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deFile", "", "", $deFile) | Out-Null
